$d = $word.ActiveDocument

# 1) Fix typo: "Not what access controls handles are granted." -> "Note what
#    access controls handles are granted."
# The "granted." run carries a grammar proofing mark (w:proofErr), so rather
# than replacing the whole sentence (which would orphan that mark), just
# retype the leading "Not" -> "Note" in place.
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext.StartsWith("Not what access controls handles are granted.")) {
        $r = $d.Range($p.Range.Start, $p.Range.Start + 3)
        $r.Text = "Note"
        break
    }
}

# 2) Remove the obsolete todo item "Add non-throwing swap." (whole bullet).
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd("`r", "`a")
    if ($ptext -eq "Add non-throwing swap.") {
        $p.Range.Delete()
        break
    }
}
